$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the content of the "AKBANK" (column D) benchmark cells that were
# removed/blanked in this update, plus the two FINASNBANK (column K)
# cells for rows 24 and 25.
$cellsToClear = @(
    "D3", "D4", "D5", "D6",
    "D8", "D9", "D10", "D11", "D12", "D13", "D14",
    "K24", "K25"
)

foreach ($cellRef in $cellsToClear) {
    $ws.Range($cellRef).ClearContents()
}
